# Update NATMI ligand-receptor pair stats (F2-Gp9 sheet) with refreshed TPM-derived values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3845463333333334
$ws.Range("H2").Value = 1.153639
$ws.Range("I2").Value = 0.1984850200147207
$ws.Range("J2").Value = 0.1984850200147207
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.20571
$ws.Range("N2").Value = 0.61713
$ws.Range("O2").Value = 0.1114078352132379
$ws.Range("P2").Value = 0.1114078352132379
$ws.Range("Q2").Value = 0.07910502623
$ws.Range("R2").Value = 0.71194523607
$ws.Range("S2").Value = 0.02211278640209623
$ws.Range("T2").Value = 0.02211278640209622
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3845463333333334
$ws.Range("H3").Value = 1.153639
$ws.Range("I3").Value = 0.1984850200147207
$ws.Range("J3").Value = 0.1984850200147207
$ws.Range("O3").Value = 0.8502405595430678
$ws.Range("P3").Value = 0.8502405595430678
$ws.Range("Q3").Value = 0.6037124914574445
$ws.Range("R3").Value = 5.433412423117001
$ws.Range("S3").Value = 0.1687600144782332
$ws.Range("T3").Value = 0.1687600144782331
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3845463333333334
$ws.Range("H4").Value = 1.153639
$ws.Range("I4").Value = 0.1984850200147207
$ws.Range("J4").Value = 0.1984850200147207
$ws.Range("M4").Value = 0.07081466666666668
$ws.Range("N4").Value = 0.212444
$ws.Range("O4").Value = 0.03835160524369437
$ws.Range("P4").Value = 0.03835160524369437
$ws.Range("Q4").Value = 0.0272315204128889
$ws.Range("R4").Value = 0.245083683716
$ws.Range("S4").Value = 0.007612219134391347
$ws.Range("T4").Value = 0.007612219134391345
$ws.Range("I5").Value = 0.5733580031870772
$ws.Range("J5").Value = 0.5733580031870772
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.20571
$ws.Range("N5").Value = 0.61713
$ws.Range("O5").Value = 0.1114078352132379
$ws.Range("P5").Value = 0.1114078352132379
$ws.Range("Q5").Value = 0.2285084278799999
$ws.Range("R5").Value = 2.05657585092
$ws.Range("S5").Value = 0.063876573937257
$ws.Range("T5").Value = 0.06387657393725701
$ws.Range("I6").Value = 0.5733580031870772
$ws.Range("J6").Value = 0.5733580031870772
$ws.Range("O6").Value = 0.8502405595430678
$ws.Range("P6").Value = 0.8502405595430678
$ws.Range("S6").Value = 0.4874922294482765
$ws.Range("T6").Value = 0.4874922294482765
$ws.Range("I7").Value = 0.5733580031870772
$ws.Range("J7").Value = 0.5733580031870772
$ws.Range("M7").Value = 0.07081466666666668
$ws.Range("N7").Value = 0.212444
$ws.Range("O7").Value = 0.03835160524369437
$ws.Range("P7").Value = 0.03835160524369437
$ws.Range("Q7").Value = 0.07866291454400001
$ws.Range("R7").Value = 0.7079662308960001
$ws.Range("S7").Value = 0.02198919980154364
$ws.Range("T7").Value = 0.02198919980154364
$ws.Range("G8").Value = 0.442033
$ws.Range("H8").Value = 1.326099
$ws.Range("I8").Value = 0.2281569767982021
$ws.Range("J8").Value = 0.2281569767982021
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.20571
$ws.Range("N8").Value = 0.61713
$ws.Range("O8").Value = 0.1114078352132379
$ws.Range("P8").Value = 0.1114078352132379
$ws.Range("Q8").Value = 0.09093060842999998
$ws.Range("R8").Value = 0.8183754758699999
$ws.Range("S8").Value = 0.02541847487388463
$ws.Range("T8").Value = 0.02541847487388463
$ws.Range("G9").Value = 0.442033
$ws.Range("H9").Value = 1.326099
$ws.Range("I9").Value = 0.2281569767982021
$ws.Range("J9").Value = 0.2281569767982021
$ws.Range("O9").Value = 0.8502405595430678
$ws.Range("P9").Value = 0.8502405595430678
$ws.Range("Q9").Value = 0.6939627831663333
$ws.Range("R9").Value = 6.245665048497
$ws.Range("S9").Value = 0.1939883156165581
$ws.Range("T9").Value = 0.1939883156165581
$ws.Range("G10").Value = 0.442033
$ws.Range("H10").Value = 1.326099
$ws.Range("I10").Value = 0.2281569767982021
$ws.Range("J10").Value = 0.2281569767982021
$ws.Range("M10").Value = 0.07081466666666668
$ws.Range("N10").Value = 0.212444
$ws.Range("O10").Value = 0.03835160524369437
$ws.Range("P10").Value = 0.03835160524369437
$ws.Range("Q10").Value = 0.03130241955066667
$ws.Range("R10").Value = 0.281721775956
$ws.Range("S10").Value = 0.008750186307759384
$ws.Range("T10").Value = 0.008750186307759384
